$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("G1").Value = "Area"
$ws.Range("H1").Value = "Atotal"
$ws.Range("J1").Value = "Atotal"
$ws.Range("K1").Value = "Qtotal"

# Row 2 - totals / summary
$ws.Range("G2").Formula = "=(D2-0)*B2/100"
$ws.Range("H2").Formula = "=SUM(G2:G11)"
$ws.Range("J2").Formula = "=H2"
$ws.Range("K2").Formula = "=F2"

# Area increments G3:G15, following the same pattern as column D/E
$ws.Range("G3").Formula = "=(D3-D2)*B3/100"

for ($r = 4; $r -le 15; $r++) {
    $prev = $r - 1
    $ws.Cells.Item($r, 7).Formula = "=(D$r-D$prev)*B$r/100"
}

$ws.Range("J2:K2").Select() | Out-Null
